$wb = $excel.ActiveWorkbook

# --- Rename "models" -> "mods" and update its selection/active state ---
$mods = $wb.Worksheets.Item("models")
$mods.Name = "mods"
$mods.Range("A35").Select() | Out-Null

# --- Add a new worksheet "covmod-comp" at the end of the workbook ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$comp = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$comp.Name = "covmod-comp"

# Column A width
$comp.Columns.Item(1).ColumnWidth = 13.14

# --- Row-label cells entered first (so shared strings land at indices 69-71) ---
$comp.Range("A2").Value = "prey"
$comp.Range("A4").Value = "closed"
$comp.Range("A5").Value = "prey*closed"

# --- Row 1 headers (entered in this order to match shared-string allocation) ---
$comp.Range("I1").Value = "MCPall"
$comp.Range("B1").Value = "KDEalla"
$comp.Range("C1").Value = "KDEallb"
$comp.Range("E1").Value = "BVb"
$comp.Range("G1").Value = "RDb"
$comp.Range("H1").Value = "RDc"
$comp.Range("D1").Value = "BVa"
$comp.Range("F1").Value = "RDa"

# --- Row 2: prey ---
$comp.Range("B2").Value = 1.12
$comp.Range("C2").Value = 1.12
$comp.Range("D2").Value = 1.3
$comp.Range("E2").Value = 1.28
$comp.Range("F2").Value = 0.78
$comp.Range("G2").Value = 0.55
$comp.Range("H2").Value = "NA"
$comp.Range("I2").Value = 0.7

# --- Row 3: human ---
$comp.Range("A3").Value = "human"
$comp.Range("B3").Value = -0.002
$comp.Range("C3").Value = -0.003
$comp.Range("D3").Value = -0.002
$comp.Range("E3").Value = -0.004
$comp.Range("F3").Value = -0.001
$comp.Range("G3").Value = -0.003
$comp.Range("H3").Value = -0.003
$comp.Range("I3").Value = -0.003

# --- Row 4: closed ---
$comp.Range("B4").Value = "NA"
$comp.Range("C4").Value = -0.75
$comp.Range("D4").Value = "NA"
$comp.Range("E4").Value = -1.01
$comp.Range("F4").Value = "NA"
$comp.Range("G4").Value = -0.1
$comp.Range("H4").Value = 0.14
$comp.Range("I4").Value = -0.72

# --- Row 5: prey*closed ---
$comp.Range("B5").Value = "NA"
$comp.Range("C5").Value = "NA"
$comp.Range("D5").Value = "NA"
$comp.Range("E5").Value = "NA"
$comp.Range("F5").Value = "NA"
$comp.Range("G5").Value = "NA"
$comp.Range("H5").Value = "NA"
$comp.Range("I5").Value = "NA"

# --- Row 6: human*closed ---
$comp.Range("A6").Value = "human*closed"
$comp.Range("B6").Value = "NA"
$comp.Range("C6").Value = 0.002
$comp.Range("D6").Value = "NA"
$comp.Range("E6").Value = 0.002
$comp.Range("F6").Value = "NA"
$comp.Range("G6").Value = 0.002
$comp.Range("H6").Value = 0.003
$comp.Range("I6").Value = 0.002

$comp.Range("A2").Select() | Out-Null
